$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 26 ("RM 232") entirely; remaining rows shift up by one.
$ws.Rows("26").Delete()

# The old row 28 ("SC 92") is now row 27 after the first shift; delete it too.
$ws.Rows("27").Delete()
